$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "29.051.07"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -0.23%  "
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.834.64"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -0.05%  "
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "0.9984"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.19%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "242.04"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.34%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "0.6118"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -3.53%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.11%  "
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.07472"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -0.86%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.2919"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.98%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "23.11"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.28%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.07682"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.43%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "1.847.50"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.73%  "
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "5.001"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.11%  "
$ws.Cells.Item(14, 5).Value = "  -0.17%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "82.51"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -1.00%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "0.000009182"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -4.52%  "
$ws.Cells.Item(17, 5).Value = "  -2.98%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "29.062.21"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.28%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "2.092.76"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.78%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "229.96"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.17%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "12.65"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.19%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -0.03%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "7.199"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.19%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.13%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "159.04"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.87%  "
$ws.Cells.Item(26, 2).Value = "Stellar"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "0.1389"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -2.02%  "
$ws.Cells.Item(27, 2).Value = "Cosmos"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "8.489"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.85%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "17.77"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -1.12%  "
$ws.Cells.Item(29, 5).Value = "  -0.59%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "4.152"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -0.06%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "4.126"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +1.22%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "0.05528"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +2.06%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "1.200"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.08%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "1.842"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.7403"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.86%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "1.141"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.10%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "2.658"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.01%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "2.769"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.26%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.01778"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -1.24%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "1.210.56"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -2.98%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "6.477"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -2.90%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.8917"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -1.35%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "1.0000"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.24%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "102.07"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.44%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "1.987.44"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.44%  "
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "65.54"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.34%  "
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.00000000123"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -1.11%  "
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "0.5086"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.68%  "
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "0.4065"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -0.21%  "
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "9.131"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +1.40%  "
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.05829"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.76%  "
